# Auto-generated script applying scheduled market-data refresh to Raiden_Profits workbook.
# For each touched cell we either set a new numeric value, or clear the cell entirely
# (when the source diff shows the cell being removed outright).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 254.85715
$ws.Range("I12").Value = 254.85715
$ws.Range("K12").Value = 254.85715
$ws.Range("M12").Value = -84.85714999999999
$ws.Range("H18").Value = 1241
$ws.Range("I18").Value = 1241
$ws.Range("K18").Value = 1241
$ws.Range("M18").Value = -957
$ws.Range("H28").Value = 10774.895
$ws.Range("I28").Value = 2508.8572
$ws.Range("K28").Value = 2508.8572
$ws.Range("M28").Value = -2023.8572
$ws.Range("H38").Value = 255.8
$ws.Range("J38").Value = 1000
$ws.Range("L38").Value = 3000
$ws.Range("N38").Value = -3744
$ws.Range("H39").Value = 641.6
$ws.Range("I39").Value = 902.6667
$ws.Range("K39").Value = 2708.0001
$ws.Range("M39").Value = -2412.0001
$ws.Range("H40").Value = 3465.5557
$ws.Range("J40").Value = 3881.6667
$ws.Range("L40").Value = 3881.6667
$ws.Range("N40").Value = -4231.6667
$ws.Range("H41").Value = 711.5
$ws.Range("I41").Value = 731.3570999999999
$ws.Range("K41").Value = 731.3570999999999
$ws.Range("M41").Value = -291.3570999999999
$ws.Range("H42").Value = 500
$ws.Range("I42").Value = 500
$ws.Range("K42").Value = 1500
$ws.Range("M42").Value = -1270
$ws.Range("H69").Value = 8745.091
$ws.Range("I69").Value = 6866
$ws.Range("J69").Value = 11000
$ws.Range("K69").Value = 20598
$ws.Range("L69").Value = 33000
$ws.Range("M69").Value = -19724
$ws.Range("N69").Value = -34748
$ws.Range("H72").Value = 8745.091
$ws.Range("I72").Value = 6866
$ws.Range("J72").Value = 11000
$ws.Range("K72").Value = 61794
$ws.Range("L72").Value = 99000
$ws.Range("M72").Value = -57426
$ws.Range("N72").Value = -107736
$ws.Range("H80").Value = 759.7273
$ws.Range("I80").Value = 1096.3334
$ws.Range("J80").Value = 633.5
$ws.Range("K80").Value = 3289.0002
$ws.Range("L80").Value = 1900.5
$ws.Range("M80").Value = -2291.0002
$ws.Range("N80").Value = -3896.5
$ws.Range("H83").Value = 759.7273
$ws.Range("I83").Value = 1096.3334
$ws.Range("J83").Value = 633.5
$ws.Range("K83").Value = 9867.000599999999
$ws.Range("L83").Value = 5701.5
$ws.Range("M83").Value = -4875.000599999999
$ws.Range("N83").Value = -15685.5
$ws.Range("H86").Value = 8955.5
$ws.Range("I86").Value = 8997.125
$ws.Range("J86").Value = 8900
$ws.Range("K86").Value = 8997.125
$ws.Range("L86").Value = 8900
$ws.Range("M86").Value = -7874.125
$ws.Range("N86").Value = -11146
$ws.Range("H89").Value = 8955.5
$ws.Range("I89").Value = 8997.125
$ws.Range("J89").Value = 8900
$ws.Range("K89").Value = 44985.625
$ws.Range("L89").Value = 44500
$ws.Range("M89").Value = -39369.625
$ws.Range("N89").Value = -55732
$ws.Range("H106").Value = 4865.778
$ws.Range("I106").Value = 4911.5
$ws.Range("J106").Value = 4500
$ws.Range("K106").Value = 4911.5
$ws.Range("L106").Value = 4500
$ws.Range("M106").Value = -4280.5
$ws.Range("N106").Value = -5762
$ws.Range("H111").Value = 638.5833
$ws.Range("I111").Value = 366.8
$ws.Range("J111").Value = 1997.5
$ws.Range("K111").Value = 1100.4
$ws.Range("L111").Value = 5992.5
$ws.Range("M111").Value = 1966.6
$ws.Range("N111").Value = -12126.5
$ws.Range("H115").Value = 1225.6666
$ws.Range("I115").Value = 1225.6666
$ws.Range("K115").Value = 3676.9998
$ws.Range("M115").Value = -2109.9998
$ws.Range("H116").Value = 4859.857
$ws.Range("I116").Value = 4853.364
$ws.Range("K116").Value = 4853.364
$ws.Range("M116").Value = -1411.364
$ws.Range("H121").Value = 999
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -6491
$ws.Range("H127").Value = 1550.6666
$ws.Range("I127").Value = 1550.6666
$ws.Range("K127").Value = 4651.9998
$ws.Range("M127").Value = 308.0002000000004
$ws.Range("H131").Value = 1595.1666
$ws.Range("I131").Value = 1595.1666
$ws.Range("K131").Value = 4785.4998
$ws.Range("M131").Value = 254.5002000000004
$ws.Range("H132").Value = 2066
$ws.Range("I132").Value = 1806.7742
$ws.Range("J132").Value = 4075
$ws.Range("K132").Value = 5420.3226
$ws.Range("L132").Value = 12225
$ws.Range("M132").Value = -2890.3226
$ws.Range("N132").Value = -17285
$ws.Range("H133").Value = 83749.5
$ws.Range("J133").Value = 83749.5
$ws.Range("L133").Value = 83749.5
$ws.Range("N133").Value = -93869.5
$ws.Range("H137").Value = 1156.75
$ws.Range("I137").Value = 1133.091
$ws.Range("J137").Value = 1208.8
$ws.Range("K137").Value = 3399.273
$ws.Range("L137").Value = 3626.4
$ws.Range("M137").Value = -849.2729999999997
$ws.Range("N137").Value = -8726.4
$ws.Range("H138").Value = 776424.8
$ws.Range("J138").Value = 855481.75
$ws.Range("L138").Value = 2566445.25
$ws.Range("N138").Value = -2576725.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2591.6667
$ws.Range("I32").Value = 1792.3433
$ws.Range("J32").Value = 13302.6
$ws.Range("K32").Value = 1792.3433
$ws.Range("L32").Value = 13302.6
$ws.Range("M32").Value = -1505.3433
$ws.Range("N32").Value = -13876.6
$ws.Range("H34").Value = 6280736.5
$ws.Range("I34").Value = 14311285
$ws.Range("J34").Value = 34754.445
$ws.Range("K34").Value = 14311285
$ws.Range("L34").Value = 34754.445
$ws.Range("M34").Value = -14311014
$ws.Range("N34").Value = -35296.445
$ws.Range("H37").Value = 17527.25
$ws.Range("J37").Value = 20025
$ws.Range("L37").Value = 20025
$ws.Range("N37").Value = -20571
$ws.Range("H63").Value = 4805.3
$ws.Range("I63").Value = 3756.625
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 3756.625
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -3070.625
$ws.Range("N63").Value = -10372
$ws.Range("H64").Value = 50091
$ws.Range("J64").Value = 50091
$ws.Range("L64").Value = 50091
$ws.Range("N64").Value = -50587
$ws.Range("H66").Value = 4805.3
$ws.Range("I66").Value = 3756.625
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 18783.125
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -15351.125
$ws.Range("N66").Value = -51864
$ws.Range("H67").Value = 50091
$ws.Range("J67").Value = 50091
$ws.Range("L67").Value = 50091
$ws.Range("N67").Value = -51807
$ws.Range("H74").Value = 2182.6
$ws.Range("I74").Value = 2182.1538
$ws.Range("K74").Value = 2182.1538
$ws.Range("M74").Value = -1308.1538
$ws.Range("H75").Value = 80000.5
$ws.Range("J75").Value = 80000.5
$ws.Range("L75").Value = 80000.5
$ws.Range("N75").Value = -81748.5
$ws.Range("H77").Value = 2182.6
$ws.Range("I77").Value = 2182.1538
$ws.Range("K77").Value = 10910.769
$ws.Range("M77").Value = -6542.769
$ws.Range("H78").Value = 80000.5
$ws.Range("J78").Value = 80000.5
$ws.Range("L78").Value = 240001.5
$ws.Range("N78").Value = -248737.5
$ws.Range("H122").Value = 5555
$ws.Range("I122").Value = 5555
$ws.Range("K122").Value = 16665
$ws.Range("M122").Value = -14215
$ws.Range("H124").Value = 57437.2
$ws.Range("J124").Value = 57437.2
$ws.Range("L124").Value = 57437.2
$ws.Range("N124").Value = -67257.2
$ws.Range("H125").Value = 66904.664
$ws.Range("J125").Value = 66904.664
$ws.Range("L125").Value = 66904.664
$ws.Range("N125").Value = -76744.664
$ws.Range("H132").Value = 3354.5
$ws.Range("I132").Value = 3048.7144
$ws.Range("J132").Value = 3782.6
$ws.Range("K132").Value = 9146.143199999999
$ws.Range("L132").Value = 11347.8
$ws.Range("M132").Value = -6616.143199999999
$ws.Range("N132").Value = -16407.8

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1539.56
$ws.Range("I20").Value = 1056.8948
$ws.Range("K20").Value = 1056.8948
$ws.Range("M20").Value = -809.8948
$ws.Range("H22").Value = 1142.25
$ws.Range("I22").Value = 471
$ws.Range("J22").Value = 4498.5
$ws.Range("K22").Value = 471
$ws.Range("L22").Value = 4498.5
$ws.Range("M22").Value = -298
$ws.Range("N22").Value = -4844.5
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H64").Value = 863.8182
$ws.Range("I64").Value = 857.625
$ws.Range("K64").Value = 857.625
$ws.Range("M64").Value = -632.625
$ws.Range("H67").Value = 863.8182
$ws.Range("I67").Value = 857.625
$ws.Range("K67").Value = 857.625
$ws.Range("M67").Value = -77.625
$ws.Range("H86").Value = 3877.652
$ws.Range("I86").Value = 4148.3887
$ws.Range("K86").Value = 4148.3887
$ws.Range("M86").Value = -3025.3887
$ws.Range("H89").Value = 3877.652
$ws.Range("I89").Value = 4148.3887
$ws.Range("K89").Value = 20741.9435
$ws.Range("M89").Value = -15125.9435
$ws.Range("H99").Value = 3840.4211
$ws.Range("I99").Value = 3678.8667
$ws.Range("J99").Value = 4446.25
$ws.Range("K99").Value = 3678.8667
$ws.Range("L99").Value = 4446.25
$ws.Range("M99").Value = -2180.8667
$ws.Range("N99").Value = -7442.25
$ws.Range("H107").Value = 2626.875
$ws.Range("I107").Value = 2814.2856
$ws.Range("J107").Value = 2364.5
$ws.Range("K107").Value = 2814.2856
$ws.Range("L107").Value = 2364.5
$ws.Range("M107").Value = -894.2856000000002
$ws.Range("N107").Value = -6204.5
$ws.Range("H132").Value = 67500
$ws.Range("J132").Value = 67500
$ws.Range("L132").Value = 67500
$ws.Range("N132").Value = -77620
$ws.Range("H134").Value = 3109.25
$ws.Range("I134").Value = 2999.4707
$ws.Range("J134").Value = 3731.3333
$ws.Range("K134").Value = 8998.4121
$ws.Range("L134").Value = 11193.9999
$ws.Range("M134").Value = -6463.4121
$ws.Range("N134").Value = -16263.9999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8999
$ws.Range("J4").Value = 8999
$ws.Range("L4").Value = 8999
$ws.Range("N4").Value = -9223
$ws.Range("H22").Value = 1480.25
$ws.Range("I22").Value = 1441.3334
$ws.Range("K22").Value = 1441.3334
$ws.Range("M22").Value = -1091.3334
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10574
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H58").Value = 3182.52
$ws.Range("I58").Value = 2098.4119
$ws.Range("K58").Value = 2098.4119
$ws.Range("M58").Value = -1895.4119
$ws.Range("H99").Value = 18194.227
$ws.Range("I99").Value = 16482
$ws.Range("J99").Value = 20248.9
$ws.Range("K99").Value = 16482
$ws.Range("L99").Value = 20248.9
$ws.Range("M99").Value = -14984
$ws.Range("N99").Value = -23244.9
$ws.Range("H107").Value = 807.94116
$ws.Range("I107").Value = 774.6667
$ws.Range("J107").Value = 887.8
$ws.Range("K107").Value = 774.6667
$ws.Range("L107").Value = 887.8
$ws.Range("M107").Value = 1145.3333
$ws.Range("N107").Value = -4727.8
$ws.Range("H124").Value = 56333
$ws.Range("J124").Value = 56333
$ws.Range("L124").Value = 56333
$ws.Range("N124").Value = -61243
$ws.Range("H126").Value = 18194.227
$ws.Range("I126").Value = 16482
$ws.Range("J126").Value = 20248.9
$ws.Range("K126").Value = 49446
$ws.Range("L126").Value = 60746.7
$ws.Range("M126").Value = -46976
$ws.Range("N126").Value = -65686.70000000001
$ws.Range("H132").Value = 2493.9614
$ws.Range("I132").Value = 2437.2917
$ws.Range("J132").Value = 3174
$ws.Range("K132").Value = 7311.875100000001
$ws.Range("L132").Value = 9522
$ws.Range("M132").Value = -4781.875100000001
$ws.Range("N132").Value = -14582
$ws.Range("H134").Value = 2469.0715
$ws.Range("I134").Value = 3740.2856
$ws.Range("J134").Value = 1197.8572
$ws.Range("K134").Value = 11220.8568
$ws.Range("L134").Value = 3593.5716
$ws.Range("M134").Value = -8685.856800000001
$ws.Range("N134").Value = -8663.571599999999
$ws.Range("H136").Value = 3182.52
$ws.Range("I136").Value = 2098.4119
$ws.Range("K136").Value = 6295.2357
$ws.Range("M136").Value = -3745.2357

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8000
$ws.Range("I3").Value = 8000
$ws.Range("K3").Value = 24000
$ws.Range("M3").Value = -23888
$ws.Range("H18").Value = 224.33333
$ws.Range("I18").Value = 224.33333
$ws.Range("K18").Value = 672.99999
$ws.Range("M18").Value = -503.99999
$ws.Range("H40").Value = 1134.3846
$ws.Range("I40").Value = 14.7
$ws.Range("K40").Value = 58.8
$ws.Range("M40").Value = 10.2
$ws.Range("H60").Value = 268.33334
$ws.Range("I60").Value = 122
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 366
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -115
$ws.Range("N60").Value = -3502
$ws.Range("H68").Value = 2245.7693
$ws.Range("J68").Value = 2561.4524
$ws.Range("L68").Value = 7684.3572
$ws.Range("N68").Value = -9306.3572
$ws.Range("H71").Value = 2245.7693
$ws.Range("J71").Value = 2561.4524
$ws.Range("L71").Value = 23053.0716
$ws.Range("N71").Value = -31165.0716
$ws.Range("H113").Value = 711907.3
$ws.Range("I113").Value = 2843477.8
$ws.Range("J113").Value = 1383.8334
$ws.Range("K113").Value = 8530433.399999999
$ws.Range("L113").Value = 4151.5002
$ws.Range("M113").Value = -8528263.399999999
$ws.Range("N113").Value = -8491.5002
$ws.Range("H129").Value = 6946.4165
$ws.Range("J129").Value = 11426.143
$ws.Range("L129").Value = 34278.429
$ws.Range("N129").Value = -44278.429
$ws.Range("H130").Value = 3499.6667
$ws.Range("I130").Value = 2999.5
$ws.Range("K130").Value = 8998.5
$ws.Range("M130").Value = -3978.5
$ws.Range("H131").Value = 17962.426
$ws.Range("I131").Value = 372203.66
$ws.Range("K131").Value = 1116610.98
$ws.Range("M131").Value = -1111570.98

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 45051.5
$ws.Range("I10").Value = 45051.5
$ws.Range("K10").Value = 45051.5
$ws.Range("M10").Value = -44882.5
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H20").Value = 49999
$ws.Range("I20").Value = 49999
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 49999
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -49754
$ws.Range("N20").ClearContents()
$ws.Range("H80").Value = 3888.4783
$ws.Range("I80").Value = 2812.6428
$ws.Range("J80").Value = 5562
$ws.Range("K80").Value = 2812.6428
$ws.Range("L80").Value = 5562
$ws.Range("M80").Value = -1814.6428
$ws.Range("N80").Value = -7558
$ws.Range("H83").Value = 3888.4783
$ws.Range("I83").Value = 2812.6428
$ws.Range("J83").Value = 5562
$ws.Range("K83").Value = 14063.214
$ws.Range("L83").Value = 27810
$ws.Range("M83").Value = -9071.214
$ws.Range("N83").Value = -37794
$ws.Range("H132").Value = 3364.7368
$ws.Range("I132").Value = 1905.7142
$ws.Range("J132").Value = 7450
$ws.Range("K132").Value = 5717.142599999999
$ws.Range("L132").Value = 22350
$ws.Range("M132").Value = -3187.142599999999
$ws.Range("N132").Value = -27410

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2162.25
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 2162.25
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214
$ws.Range("H40").Value = 1895
$ws.Range("I40").Value = 1895
$ws.Range("K40").Value = 1895
$ws.Range("M40").Value = -1759
$ws.Range("H46").Value = 3013.2
$ws.Range("I46").Value = 1731.3334
$ws.Range("J46").Value = 3562.5715
$ws.Range("K46").Value = 1731.3334
$ws.Range("L46").Value = 3562.5715
$ws.Range("M46").Value = -1543.3334
$ws.Range("N46").Value = -3938.5715
$ws.Range("H55").Value = 276.14285
$ws.Range("I55").Value = 271.83334
$ws.Range("J55").Value = 302
$ws.Range("K55").Value = 271.83334
$ws.Range("L55").Value = 302
$ws.Range("M55").Value = -98.83334000000002
$ws.Range("N55").Value = -648
$ws.Range("H68").Value = 2928.8572
$ws.Range("I68").Value = 2833.6667
$ws.Range("K68").Value = 2833.6667
$ws.Range("M68").Value = -2084.6667
$ws.Range("H71").Value = 2928.8572
$ws.Range("I71").Value = 2833.6667
$ws.Range("K71").Value = 14168.3335
$ws.Range("M71").Value = -10424.3335
$ws.Range("H104").Value = 35857.57
$ws.Range("J104").Value = 36834
$ws.Range("L104").Value = 36834
$ws.Range("N104").Value = -43822
$ws.Range("H122").Value = 2849.4
$ws.Range("I122").Value = 2849.4
$ws.Range("K122").Value = 8548.200000000001
$ws.Range("M122").Value = -6098.200000000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H18").Value = 5916.6665
$ws.Range("I18").Value = 4000
$ws.Range("J18").Value = 7833.3335
$ws.Range("K18").Value = 4000
$ws.Range("L18").Value = 7833.3335
$ws.Range("M18").Value = -3827
$ws.Range("N18").Value = -8179.3335
$ws.Range("H22").Value = 16670.666
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 25000
$ws.Range("K22").Value = 12
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 281
$ws.Range("N22").Value = -25586
$ws.Range("H31").Value = 34444
$ws.Range("I31").Value = 34444
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 34444
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -34096
$ws.Range("N31").ClearContents()
$ws.Range("H52").Value = 55428.57
$ws.Range("I52").Value = 28000
$ws.Range("J52").Value = 60000
$ws.Range("K52").Value = 28000
$ws.Range("L52").Value = 60000
$ws.Range("M52").Value = -27774
$ws.Range("N52").Value = -60452
$ws.Range("H62").Value = 9053
$ws.Range("J62").Value = 9719
$ws.Range("L62").Value = 9719
$ws.Range("N62").Value = -10967
$ws.Range("H65").Value = 9053
$ws.Range("J65").Value = 9719
$ws.Range("L65").Value = 48595
$ws.Range("N65").Value = -54835
$ws.Range("H97").Value = 38464
$ws.Range("J97").Value = 38464
$ws.Range("L97").Value = 38464
$ws.Range("N97").Value = -40446
$ws.Range("H107").Value = 1273.6
$ws.Range("I107").Value = 1267.4
$ws.Range("J107").Value = 1279.8
$ws.Range("K107").Value = 3802.2
$ws.Range("L107").Value = 3839.4
$ws.Range("M107").Value = -1882.2
$ws.Range("N107").Value = -7679.4
$ws.Range("H113").Value = 1370
$ws.Range("I113").Value = 1333.7142
$ws.Range("J113").Value = 1497
$ws.Range("K113").Value = 4001.1426
$ws.Range("L113").Value = 4491
$ws.Range("M113").Value = -1831.1426
$ws.Range("N113").Value = -8831
$ws.Range("H121").Value = 59998
$ws.Range("J121").Value = 59998
$ws.Range("L121").Value = 59998
$ws.Range("N121").Value = -63492
$ws.Range("H122").Value = 4334.2354
$ws.Range("J122").Value = 2198
$ws.Range("L122").Value = 6594
$ws.Range("N122").Value = -11494
$ws.Range("H132").Value = 5172.343
$ws.Range("I132").Value = 3808.5356
$ws.Range("K132").Value = 11425.6068
$ws.Range("M132").Value = -8895.606800000001
$ws.Range("H136").Value = 3888.875
$ws.Range("I136").Value = 4864.5
$ws.Range("J136").Value = 1937.625
$ws.Range("K136").Value = 14593.5
$ws.Range("L136").Value = 5812.875
$ws.Range("M136").Value = -12043.5
$ws.Range("N136").Value = -10912.875
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
